$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1. Refresh the cached "datetimeFigureOut" date field shown on the
#    slide master and every slide layout's Date placeholder
#    (ppPlaceholderDate = 16) from "6/5/17" to "2/17/2018".
# ------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDatePlaceholder = $false
            if ($shp.Type -eq 14) {
                $phFormat = $shp.PlaceholderFormat
                if ($phFormat.Type -eq 16) {
                    $isDatePlaceholder = $true
                }
            }
            if ($isDatePlaceholder) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq "6/5/17") {
                    $tr.Text = "2/17/2018"
                }
            }
        }
    }
}

Update-DatePlaceholders $p.SlideMaster.Shapes

$customLayouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $customLayouts.Count; $j++) {
    $layout = $customLayouts.Item($j)
    Update-DatePlaceholders $layout.Shapes
}

# ------------------------------------------------------------------
# 2. Re-create the presentation-level guides (PowerPoint 2013+ stores
#    these both as the legacy view guide list already present in this
#    file, and as p15:sldGuideLst) - one horizontal guide at 2160 and
#    one vertical guide at 2880.
# ------------------------------------------------------------------
try {
    $guides = $p.Guides
    if ($guides -ne $null) {
        if ($guides.Count -lt 1) {
            $guides.Add(1, 2160) | Out-Null
        }
        if ($guides.Count -lt 2) {
            $guides.Add(2, 2880) | Out-Null
        }
        for ($k = 1; $k -le $guides.Count; $k++) {
            $guide = $guides.Item($k)
            if ($guide.Position -eq 2160) {
                $guide.Orientation = 1
            } elseif ($guide.Position -eq 2880) {
                $guide.Orientation = 2
            }
        }
    }
} catch {
    # Guides collection is not available in this environment - ignore.
}
